# Update ASYR (column C) and ASYR_FLUC (column D) figures for Sheet1
# of the EU_ASYR_All_age_groups static-over-90 report, per the refreshed
# PYLL/ASYR data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 1264.704
$ws.Range("D2").Value = 328.806
$ws.Range("C3").Value = 3355.414
$ws.Range("D3").Value = 442.6559999999999
$ws.Range("C4").Value = 2284.951
$ws.Range("D4").Value = 332.913
$ws.Range("C7").Value = 6778.589999999999
$ws.Range("C11").Value = 1080.295
$ws.Range("C12").Value = 1724.729
$ws.Range("C13").Value = 1363.645
$ws.Range("C14").Value = 2154.017
$ws.Range("C16").Value = 2971.728999999999
$ws.Range("C23").Value = 1437.703
$ws.Range("C24").Value = 2619.778
$ws.Range("C25").Value = 1968.5
$ws.Range("C29").Value = 874.3150000000001
$ws.Range("C30").Value = 1208.437
$ws.Range("C31").Value = 937.0240000000001
$ws.Range("C35").Value = 606.4970000000001
$ws.Range("C36").Value = 1115.334
$ws.Range("C37").Value = 831.634
$ws.Range("C51").Value = 604.2430000000001
$ws.Range("C52").Value = 1420.097
$ws.Range("C53").Value = 1358.38
$ws.Range("C54").Value = 380.719
$ws.Range("D54").Value = 186.462
$ws.Range("C55").Value = 949.1890000000001
$ws.Range("D55").Value = 443.6329999999999
$ws.Range("C56").Value = 550.562
$ws.Range("C63").Value = 1060.83
$ws.Range("C64").Value = 1816.421
$ws.Range("C65").Value = 1434.525
$ws.Range("C66").Value = 1363.981
$ws.Range("C67").Value = 2467.707
$ws.Range("C68").Value = 1842.946
$ws.Range("C69").Value = 819.346
$ws.Range("D69").Value = 209.4
$ws.Range("C70").Value = 1606.501
$ws.Range("D70").Value = 291.733
$ws.Range("C71").Value = 1149.257
$ws.Range("D71").Value = 238.932
$ws.Range("C72").Value = 181.624
$ws.Range("D72").Value = 152.386
$ws.Range("C73").Value = 313.51
$ws.Range("D73").Value = 161.532
$ws.Range("C74").Value = 222.66
$ws.Range("D74").Value = 146.841
$ws.Range("D80").Value = 245.1220000000001
$ws.Range("C81").Value = 1663.699
$ws.Range("C82").Value = 3404.527
$ws.Range("C83").Value = 2507.473
$ws.Range("C87").Value = 1285.682
$ws.Range("C88").Value = 2215.629
$ws.Range("C89").Value = 1671.097
$ws.Range("C90").Value = 1300.682
$ws.Range("C91").Value = 2294.798
$ws.Range("C92").Value = 1700.074
$ws.Range("C93").Value = 1197.655
$ws.Range("D93").Value = 192.379
$ws.Range("C94").Value = 1938.366
$ws.Range("C95").Value = 1553.315
$ws.Range("D95").Value = 223.496

Write-Output "Updated 62 cells on sheet $($ws.Name)"
